$d = $word.ActiveDocument

# Apply text replacements using explicit character-offset Range assignment.
# (Find.Execute in this runtime only matches patterns that start exactly at a
#  run boundary and does not tolerate leading/trailing whitespace in the search
#  string, so direct Range offsets are used instead. Edits are applied from the
#  end of the document towards the start so earlier offsets stay valid.)

$d.Range(2253, 2439).Text = ": los apasionados del té, las personas conscientes de la salud, los amantes de las bebidas calientes especiadas, y cualquier persona que desea explorar los sabores intensos del chai indio tradicional."
$d.Range(2046, 2152).Text = ": respaldamos nuestro producto y ofrecemos una garantía de satisfacción."
$d.Range(2038, 2046).Text = "Garantía de satisfacción del cliente"
$d.Range(1855, 2037).Text = ": el té chai Mystic Spice viene en un envase elegante, ecológico, lo que lo convierte en el regalo ideal para los amantes del té o un capricho lujoso para ti mismo."
$d.Range(1844, 1855).Text = "Envase elegante"
$d.Range(1608, 1843).Text = ": al estar comprometidos con la sostenibilidad, obtenemos nuestros ingredientes de pequeñas explotaciones que practican la agricultura ecológica, garantizando no solo la mejor calidad, sino también el bienestar de nuestro planeta."
$d.Range(1602, 1608).Text = "Origen sostenible"
$d.Range(1263, 1481).Text = ": ya sea que te guste tu chai bien caliente, o prefieras un refrescante té helado, o un latte cremoso, nuestra mezcla es lo suficientemente versátil como para adaptarse a cualquier preferencia."
$d.Range(1255, 1263).Text = "Opciones versátiles de preparación"
$d.Range(992, 1144).Text = ": el aroma cálido y especiado, y el sabor profundo y vigorizante de nuestro chai hacen que sea la bebida perfecta para comenzar el día o relajarse por la noche."
$d.Range(974, 992).Text = "Aroma y sabor intensos"
$d.Range(828, 973).Text = "El jengibre y el cardamomo ayudan a la digestión, la canela ayuda a regular el azúcar en sangre y el clavo aumento los antioxidantes."
$d.Range(706, 827).Text = ": cada ingrediente del té chai premium Mystic Spice se elige por sus beneficios naturales para la salud."
$d.Range(694, 706).Text = "Ingredientes beneficiosos para la salud"
$d.Range(441, 619).Text = ": nuestro chai es una mezcla armoniosa de hojas de té negro premium y una selección de especias molidas, incluyendo canela, cardamomo, clavo, jengibre y pimienta negra."
$d.Range(435, 441).Text = "Mezcla auténtica"
$d.Range(406, 434).Text = "Características principales:"
$d.Range(6, 51).Text = ": Té chai premium Mystic Spice"
$d.Range(0, 6).Text = "Nombre de producto"

# Make the "Características principales:" label run bold (w:b val 0 -> w:b)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Características principales:") {
        $p.Range.Font.Bold = 1
    }
}
